$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 ("Sửa lại Usecase") now also marked as completed
$ws.Range("D10").Value = "Hoàn thành"

# New row content: "Đặc tả Usecase" inserted into the Tuần 2 block (between
# "Sửa lại Usecase" and "Vẽ Class Diagram"), marked as completed. Since this
# falls inside an already-merged/bordered block (B10:B13) the row beneath it
# just gets overwritten in place rather than physically inserting a row.
$ws.Range("C11").Value = "Đặc tả Usecase"
$ws.Range("D11").Value = "Hoàn thành"

# The new row pushes the remaining tasks down by one position
$ws.Range("C12").Value = "Vẽ Class Diagram"
$ws.Range("C13").Value = "Vẽ Sequence Diagram"
$ws.Range("D13").Value = "Hoàn thành"

# "Thiết kế giao diện" no longer appears anywhere in the table

# The new "Đặc tả Usecase" cell sits inside the same visual box as its
# neighbours, so it shouldn't draw an internal horizontal divider above or
# below itself.
$c11 = $ws.Range("C11")
$c11.Borders.Item(8).LineStyle = 0
$c11.Borders.Item(9).LineStyle = 0

# Leave the selection where the author last left it
$ws.Range("F6:F7").Select()
